$wb = $excel.ActiveWorkbook
$oldName = $wb.ActiveSheet.Name

# Add a new sheet and drop the old placeholder sheet so the sheetId
# counter advances the same way it did in the authored workbook
# (NONE/sheetId=3 -> USER_LIST/sheetId=4).
$ws = $wb.Worksheets.Add()
$ws.Name = "USER_LIST"
[void]$wb.Worksheets.Item($oldName).Delete()

# Write the terminator string first so it lands at shared-string index 0,
# matching the authored workbook's string table order, then fill in the
# header and data rows.
$ws.Range("A4").Value = "*END*"

# Header row
$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "User Group"

# Data rows
$ws.Range("A2").Value = "LEO"
$ws.Range("B2").Value = "CLIENT_ADMIN"

$ws.Range("A3").Value = "NEMOO"
$ws.Range("B3").Value = "CLIENT_ADMIN"

# Column widths to match bestFit sizing in the target file
# (closest achievable values given the runtime's column-width quantization)
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 12.833333333333334

# Selection matching the target sheet view
[void]$ws.Range("A4:XFD5").Select()
